$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.109.58'
$ws.Range("E2").Value = '  -1.18%  '

$ws.Range("D3").Value = '3.580.63'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  -0.07%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '236.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.19%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '650.50'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.01%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.48'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.48%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("E10").Value = '  -2.42%  '

$ws.Range("D11").Value = '3.580.46'
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("E12").Value = '  +1.04%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '42.32'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.02%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.49'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").Value = '4.240.10'
$ws.Range("E15").Value = '  -0.82%  '

$ws.Range("D16").Value = '95.023.17'
$ws.Range("E16").Value = '  -1.20%  '

$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").Value = '3.565.69'
$ws.Range("E18").Value = '  -0.11%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.74'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.56%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.52'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -5.73%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.88'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.01%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '3.46'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '508.00'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.71%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.481'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.99%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("E26").Value = '  -1.64%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '95.27'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.69%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '12.52'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.52%  '

$ws.Range("D29").Value = '3.771.44'
$ws.Range("E29").Value = '  -0.14%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '3.01'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.47%  '

$ws.Range("E31").Value = '  -1.16%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '11.47'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '

$ws.Range("E33").Value = '  +0.33%  '

$ws.Range("E34").Value = '  -0.37%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.177'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.56%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '31.80'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +4.41%  '

$ws.Range("E37").Value = '  -1.22%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.65'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +10.17%  '

$ws.Range("E39").Value = '  +7.62%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '580.31'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("E42").Value = '  -0.91%  '

$ws.Range("E43").Value = '  -2.64%  '

$ws.Range("E44").Value = '  +1.40%  '

$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '5.73'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.29'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.96%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '34.00'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +31.03%  '

$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '23.38'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.91%  '

$ws.Range("E49").Value = '  -4.42%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '3.55'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '53.34'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.84%  '
